$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 28533.191
$ws.Range("I98").Value = 31289.055
$ws.Range("J98").Value = 11998
$ws.Range("K98").Value = 31289.055
$ws.Range("L98").Value = 11998
$ws.Range("M98").Value = -29791.055
$ws.Range("N98").Value = -14994

$ws.Range("H122").Value = 28533.191
$ws.Range("I122").Value = 31289.055
$ws.Range("J122").Value = 11998
$ws.Range("K122").Value = 93867.16500000001
$ws.Range("L122").Value = 35994
$ws.Range("M122").Value = -91417.16500000001
$ws.Range("N122").Value = -40894

$ws.Range("H135").Value = 3608.2083
$ws.Range("I135").Value = 4336.8423
$ws.Range("K135").Value = 39031.58070000001
$ws.Range("M135").Value = -36496.58070000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6437.484
$ws.Range("I32").Value = 6437.484
$ws.Range("K32").Value = 6437.484
$ws.Range("M32").Value = -6150.484

$ws.Range("H98").Value = 355000
$ws.Range("J98").Value = 355000
$ws.Range("L98").Value = 355000
$ws.Range("N98").Value = -360990

$ws.Range("H122").Value = 1027157.44
$ws.Range("I122").Value = 3077.0334
$ws.Range("J122").Value = 3390419.8
$ws.Range("K122").Value = 9231.100199999999
$ws.Range("L122").Value = 10171259.4
$ws.Range("M122").Value = -6781.100199999999
$ws.Range("N122").Value = -10176159.4

$ws.Range("H132").Value = 8078.484
$ws.Range("I132").Value = 10060.833
$ws.Range("K132").Value = 30182.499
$ws.Range("M132").Value = -27652.499

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H43").Value = 249859.14
$ws.Range("J43").Value = 249859.14
$ws.Range("L43").Value = 249859.14
$ws.Range("N43").Value = -250221.14

$ws.Range("H86").Value = 4655
$ws.Range("I86").Value = 5293.864
$ws.Range("K86").Value = 5293.864
$ws.Range("M86").Value = -4170.864

$ws.Range("H89").Value = 4655
$ws.Range("I89").Value = 5293.864
$ws.Range("K89").Value = 26469.32
$ws.Range("M89").Value = -20853.32

$ws.Range("H94").Value = 9636.727999999999
$ws.Range("I94").Value = 12512.261
$ws.Range("J94").Value = 3023
$ws.Range("K94").Value = 12512.261
$ws.Range("L94").Value = 3023
$ws.Range("M94").Value = -12061.261
$ws.Range("N94").Value = -3925

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 14841.25
$ws.Range("I31").Value = 23016.666
$ws.Range("K31").Value = 23016.666
$ws.Range("M31").Value = -22721.666

$ws.Range("H34").Value = 14841.25
$ws.Range("I34").Value = 23016.666
$ws.Range("K34").Value = 23016.666
$ws.Range("M34").Value = -22814.666

$ws.Range("H41").Value = 500
$ws.Range("I41").Value = 500
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 500
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -72
$ws.Range("N41").ClearContents()

$ws.Range("H107").Value = 9272.346
$ws.Range("I107").Value = 10020.042
$ws.Range("J107").Value = 300
$ws.Range("K107").Value = 10020.042
$ws.Range("L107").Value = 300
$ws.Range("M107").Value = -8100.041999999999
$ws.Range("N107").Value = -4140

$ws.Range("H122").Value = 7418.8
$ws.Range("I122").Value = 11241.667
$ws.Range("J122").Value = 1684.5
$ws.Range("K122").Value = 33725.001
$ws.Range("L122").Value = 5053.5
$ws.Range("M122").Value = -31275.001
$ws.Range("N122").Value = -9953.5

$ws.Range("H132").Value = 17233.223
$ws.Range("J132").Value = 69697
$ws.Range("L132").Value = 209091
$ws.Range("N132").Value = -214151

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 863.3333
$ws.Range("I22").Value = 667.1429000000001
$ws.Range("J22").Value = 1550
$ws.Range("K22").Value = 2001.4287
$ws.Range("L22").Value = 4650
$ws.Range("M22").Value = -1832.4287
$ws.Range("N22").Value = -4988

$ws.Range("H27").Value = 863.3333
$ws.Range("I27").Value = 667.1429000000001
$ws.Range("J27").Value = 1550
$ws.Range("K27").Value = 2001.4287
$ws.Range("L27").Value = 4650
$ws.Range("M27").Value = -1899.4287
$ws.Range("N27").Value = -4854

$ws.Range("H63").Value = 2428.25
$ws.Range("I63").Value = 1899.5
$ws.Range("J63").Value = 2957
$ws.Range("K63").Value = 5698.5
$ws.Range("L63").Value = 8871
$ws.Range("M63").Value = -4949.5
$ws.Range("N63").Value = -10369

$ws.Range("H64").Value = 2713.5715
$ws.Range("I64").Value = 2799.2
$ws.Range("J64").Value = 2499.5
$ws.Range("K64").Value = 8397.599999999999
$ws.Range("L64").Value = 7498.5
$ws.Range("M64").Value = -8127.599999999999
$ws.Range("N64").Value = -8038.5

$ws.Range("H66").Value = 2428.25
$ws.Range("I66").Value = 1899.5
$ws.Range("J66").Value = 2957
$ws.Range("K66").Value = 17095.5
$ws.Range("L66").Value = 26613
$ws.Range("M66").Value = -13351.5
$ws.Range("N66").Value = -34101

$ws.Range("H67").Value = 2713.5715
$ws.Range("I67").Value = 2799.2
$ws.Range("J67").Value = 2499.5
$ws.Range("K67").Value = 8397.599999999999
$ws.Range("L67").Value = 7498.5
$ws.Range("M67").Value = -7461.599999999999
$ws.Range("N67").Value = -9370.5

$ws.Range("H87").Value = 16906.54
$ws.Range("I87").Value = 11510.5
$ws.Range("J87").Value = 19304.777
$ws.Range("K87").Value = 34531.5
$ws.Range("L87").Value = 57914.33099999999
$ws.Range("M87").Value = -33283.5
$ws.Range("N87").Value = -60410.33099999999

$ws.Range("H90").Value = 16906.54
$ws.Range("I90").Value = 11510.5
$ws.Range("J90").Value = 19304.777
$ws.Range("K90").Value = 103594.5
$ws.Range("L90").Value = 173742.993
$ws.Range("M90").Value = -97354.5
$ws.Range("N90").Value = -186222.993

$ws.Range("H103").Value = 1314.3158
$ws.Range("I103").Value = 696.4286
$ws.Range("J103").Value = 1674.75
$ws.Range("K103").Value = 2089.2858
$ws.Range("L103").Value = 5024.25
$ws.Range("M103").Value = -1210.2858
$ws.Range("N103").Value = -6782.25

$ws.Range("H108").Value = 2000
$ws.Range("I108").Value = 2000
$ws.Range("K108").Value = 6000
$ws.Range("M108").Value = -3120

$ws.Range("H117").Value = 3082
$ws.Range("I117").Value = 3997.3333
$ws.Range("J117").Value = 2166.6667
$ws.Range("K117").Value = 11991.9999
$ws.Range("L117").Value = 6500.000100000001
$ws.Range("M117").Value = -8549.999899999999
$ws.Range("N117").Value = -13384.0001

$ws.Range("H121").Value = 1335628.4
$ws.Range("I121").Value = 685.6
$ws.Range("J121").Value = 2003099.8
$ws.Range("K121").Value = 2056.8
$ws.Range("L121").Value = 6009299.4
$ws.Range("M121").Value = -746.8000000000002
$ws.Range("N121").Value = -6011919.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 175.4
$ws.Range("I2").Value = 154.64706
$ws.Range("J2").Value = 293
$ws.Range("K2").Value = 154.64706
$ws.Range("L2").Value = 293
$ws.Range("M2").Value = -41.64706000000001
$ws.Range("N2").Value = -519

$ws.Range("H70").Value = 7187.3335
$ws.Range("I70").Value = 6460
$ws.Range("J70").Value = 8096.5
$ws.Range("K70").Value = 6460
$ws.Range("L70").Value = 8096.5
$ws.Range("M70").Value = -6190
$ws.Range("N70").Value = -8636.5

$ws.Range("H73").Value = 7187.3335
$ws.Range("I73").Value = 6460
$ws.Range("J73").Value = 8096.5
$ws.Range("K73").Value = 6460
$ws.Range("L73").Value = 8096.5
$ws.Range("M73").Value = -5524
$ws.Range("N73").Value = -9968.5

$ws.Range("H102").Value = 4548.4146
$ws.Range("I102").Value = 5444.2334
$ws.Range("K102").Value = 5444.2334
$ws.Range("M102").Value = -3822.2334

$ws.Range("H126").Value = 11564.228
$ws.Range("I126").Value = 17666.875
$ws.Range("J126").Value = 8077
$ws.Range("K126").Value = 53000.625
$ws.Range("L126").Value = 24231
$ws.Range("M126").Value = -50530.625
$ws.Range("N126").Value = -29171

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 2000
$ws.Range("I14").Value = 2000
$ws.Range("K14").Value = 2000
$ws.Range("M14").Value = -1828

$ws.Range("H103").Value = 60199.332
$ws.Range("J103").Value = 60199.332
$ws.Range("L103").Value = 60199.332
$ws.Range("N103").Value = -62543.332

$ws.Range("H122").Value = 4111.081
$ws.Range("I122").Value = 5264.25
$ws.Range("J122").Value = 3557.56
$ws.Range("K122").Value = 15792.75
$ws.Range("L122").Value = 10672.68
$ws.Range("M122").Value = -13342.75
$ws.Range("N122").Value = -15572.68

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3835.2354
$ws.Range("I122").Value = 1615.7368
$ws.Range("K122").Value = 4847.2104
$ws.Range("M122").Value = -2397.2104

$ws.Range("H132").Value = 6605.92
$ws.Range("I132").Value = 7775.7407
$ws.Range("K132").Value = 23327.2221
$ws.Range("M132").Value = -20797.2221
